# Apply updated crypto market data to the sheet (price + 1h volume change)
# per coinranking.com snapshot refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.283.46"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").Value = "2.008.32"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.00"
$ws.Range("E5").Value = "  +4.63%  "

$ws.Range("D6").Value = "0.613"
$ws.Range("E6").Value = "  -1.55%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.00"
$ws.Range("E8").Value = "  -6.01%  "

$ws.Range("D9").Value = "0.384"
$ws.Range("E9").Value = "  -0.59%  "

$ws.Range("D10").Value = "0.0768"
$ws.Range("E10").Value = "  -5.27%  "

$ws.Range("E11").Value = "  -1.94%  "

$ws.Range("D12").Value = "2.305.06"
$ws.Range("E12").Value = "  +0.26%  "

$ws.Range("D13").Value = "14.22"
$ws.Range("E13").Value = "  -5.30%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "21.02"
$ws.Range("E14").Value = "  -5.65%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.800"
$ws.Range("E15").Value = "  -5.14%  "

$ws.Range("D16").Value = "5.23"
$ws.Range("E16").Value = "  -4.12%  "

$ws.Range("D17").Value = "2.009.04"
$ws.Range("E17").Value = "  +0.36%  "

$ws.Range("D18").Value = "37.181.61"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("D19").Value = "69.54"
$ws.Range("E19").Value = "  -0.92%  "

$ws.Range("D20").Value = "0.0₃0833"
$ws.Range("E20").Value = "  -3.75%  "

$ws.Range("D21").Value = "5.13"
$ws.Range("E21").Value = "  -1.30%  "

$ws.Range("D22").Value = "228.38"
$ws.Range("E22").Value = "  -0.85%  "

$ws.Range("D23").Value = "2.62"
$ws.Range("E23").Value = "  +5.83%  "

$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").Value = "164.68"
$ws.Range("E26").Value = "  +0.20%  "

$ws.Range("D27").Value = "8.89"
$ws.Range("E27").Value = "  -5.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.60"
$ws.Range("E28").Value = "  -0.21%  "

$ws.Range("E29").Value = "  -6.78%  "

$ws.Range("D30").Value = "1.33"
$ws.Range("E30").Value = "  -3.04%  "

$ws.Range("D31").Value = "0.119"
$ws.Range("E31").Value = "  -1.30%  "

$ws.Range("D32").Value = "4.63"
$ws.Range("E32").Value = "  -3.50%  "

$ws.Range("D33").Value = "0.0642"
$ws.Range("E33").Value = "  -1.91%  "

$ws.Range("D34").Value = "4.53"
$ws.Range("E34").Value = "  +1.69%  "

$ws.Range("D35").Value = "2.37"
$ws.Range("E35").Value = "  -2.56%  "

$ws.Range("E36").Value = "  +0.89%  "

$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "3.36"
$ws.Range("E38").Value = "  -1.56%  "

$ws.Range("D39").Value = "5.24"
$ws.Range("E39").Value = "  -1.33%  "

$ws.Range("E40").Value = "  +3.79%  "

$ws.Range("D41").Value = "1.21"
$ws.Range("E41").Value = "  +2.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0930"
$ws.Range("E42").Value = "  -5.46%  "

$ws.Range("E43").Value = "  -0.68%  "

$ws.Range("D44").Value = "1.402.25"
$ws.Range("E44").Value = "  +2.41%  "

$ws.Range("D45").Value = "89.46"
$ws.Range("E45").Value = "  -2.66%  "

$ws.Range("D46").Value = "15.68"
$ws.Range("E46").Value = "  -5.28%  "

$ws.Range("E47").Value = "  -2.51%  "

$ws.Range("D48").Value = "7.04"
$ws.Range("E48").Value = "  -4.39%  "

$ws.Range("E49").Value = "  +2.41%  "

$ws.Range("D50").Value = "2.197.92"
$ws.Range("E50").Value = "  +0.32%  "

$ws.Range("E51").Value = "  -7.44%  "
